$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.108.06"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "3.767.27"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").Value = "3.767.67"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +5.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.88"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000249"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "4.377.40"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "3.748.29"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "69.230.01"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.46"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.76"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.74"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.45"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "3.902.13"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").Value = "3.699.77"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.14"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.80%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.327"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "427.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.66"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.48"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.19"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  +10.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "2.798.38"
$ws.Range("E51").Value = "  -0.47%  "
